$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: Andrew Nembhard -> Luka Doncic
$ws.Range("A12").Value = "Luka Doncic"
$ws.Range("B12").Value = "PG,SG"
$ws.Range("C12").Value = "Dallas Mavericks"

# Row 13: Nikola Vucevic -> Jrue Holiday (new player inserted here)
$ws.Range("A13").Value = "Jrue Holiday"
$ws.Range("B13").Value = "PG,SG"
$ws.Range("C13").Value = "Boston Celtics"

# Row 14: De'Aaron Fox -> Nikola Vucevic (shifted down from old row 13)
$ws.Range("A14").Value = "Nikola Vucevic"
$ws.Range("B14").Value = "PF,C"
$ws.Range("C14").Value = "Chicago Bulls"

# Row 15: DeMar DeRozan -> De'Aaron Fox (shifted down from old row 14)
$ws.Range("A15").Value = "De'Aaron Fox"
$ws.Range("B15").Value = "PG"
$ws.Range("C15").Value = "Sacramento Kings"

# Row 16: Luka Doncic -> DeMar DeRozan (shifted down from old row 15)
$ws.Range("A16").Value = "DeMar DeRozan"
$ws.Range("B16").Value = "SF,PF"
$ws.Range("C16").Value = "Sacramento Kings"

# Row 17: Josh Giddey stays Josh Giddey (unchanged)
$ws.Range("A17").Value = "Josh Giddey"
$ws.Range("B17").Value = "PG,SG,SF"
$ws.Range("C17").Value = "Chicago Bulls"

# Row 18: new row - Goga Bitadze
$ws.Range("A18").Value = "Goga Bitadze"
$ws.Range("B18").Value = "C"
$ws.Range("C18").Value = "Orlando Magic"
